$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.342.41"
$ws.Range("D2").NumberFormat = "General"
$ws.Range("E2").Value = "  -2.55%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.687.22"
$ws.Range("D3").NumberFormat = "General"
$ws.Range("E3").Value = "  -3.34%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "691.45"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  -1.93%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "162.81"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  -5.42%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.686.06"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = "  -3.34%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  -4.50%  "
$ws.Range("E10").Value = "  -8.42%  "
$ws.Range("E11").Value = "  -3.25%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.445"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = "  -3.69%  "
$ws.Range("E13").Value = "  -5.06%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.54"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  -6.90%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.306.89"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = "  -3.38%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.684.94"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = "  -3.80%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "69.370.39"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = "  -2.52%  "
$ws.Range("E18").Value = "  -1.03%  "
$ws.Range("E19").Value = "  -6.99%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.62"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = "  -7.70%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "481.75"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  -5.83%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.99"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = "  -6.71%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.666"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = "  -7.93%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "80.07"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  -4.82%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.831.93"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  -3.39%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000131"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = "  -9.33%  "
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.43"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = "  -5.31%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.56"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "  -8.11%  "
$ws.Range("E30").Value = "  -9.98%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.73"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = "  -9.71%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.86"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  -7.49%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.09"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = "  -7.00%  "
$ws.Range("E34").Value = "  -6.98%  "
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("E36").Value = "  -4.35%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.653.02"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = "  -3.26%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.53"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = "  -7.15%  "
$ws.Range("E39").Value = "  +4.92%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.32"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = "  -2.24%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0933"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = "  -7.75%  "
$ws.Range("E43").Value = "  -0.06%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.952"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  -7.04%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "163.48"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  -4.88%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "48.03"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = "  -3.34%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.84"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  -13.99%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "30.04"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = "  +1.88%  "
$ws.Range("E49").Value = "  +0.67%  "
$ws.Range("E50").Value = "  -7.51%  "
$ws.Range("E51").Value = "  -0.51%  "
